$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 873.3333
$ws.Range("J97").Value = 873.3333
$ws.Range("L97").Value = 2619.9999
$ws.Range("N97").Value = -3611.9999

$ws.Range("H99").Value = 616.1111
$ws.Range("I99").Value = 510.66666
$ws.Range("J99").Value = 827
$ws.Range("K99").Value = 1531.99998
$ws.Range("L99").Value = 2481
$ws.Range("M99").Value = -33.99998000000005
$ws.Range("N99").Value = -5477

$ws.Range("H101").Value = 1452.7
$ws.Range("I101").Value = 385.6
$ws.Range("J101").Value = 2519.8
$ws.Range("K101").Value = 1156.8
$ws.Range("L101").Value = 7559.400000000001
$ws.Range("M101").Value = 465.1999999999998
$ws.Range("N101").Value = -10803.4

$ws.Range("H129").Value = 941.44446
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 966.11536
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 2898.34608
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -12898.34608

$ws.Range("H137").Value = 3167.7585
$ws.Range("I137").Value = 1878.8948
$ws.Range("J137").Value = 5616.6
$ws.Range("K137").Value = 5636.6844
$ws.Range("L137").Value = 16849.8
$ws.Range("M137").Value = -3086.6844
$ws.Range("N137").Value = -21949.8

$ws.Range("H138").Value = 1606.3
$ws.Range("I138").Value = 776.25
$ws.Range("J138").Value = 1813.8125
$ws.Range("K138").Value = 2328.75
$ws.Range("L138").Value = 5441.4375
$ws.Range("M138").Value = 2811.25
$ws.Range("N138").Value = -15721.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6790.367
$ws.Range("I32").Value = 5460.844
$ws.Range("K32").Value = 5460.844
$ws.Range("M32").Value = -5173.844

$ws.Range("H45").Value = 2030.125
$ws.Range("I45").Value = 1950.6666
$ws.Range("J45").Value = 2268.5
$ws.Range("K45").Value = 1950.6666
$ws.Range("L45").Value = 2268.5
$ws.Range("M45").Value = -1573.6666
$ws.Range("N45").Value = -3022.5

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H61").Value = 5118.161
$ws.Range("I61").Value = 4143.227
$ws.Range("J61").Value = 7501.3335
$ws.Range("K61").Value = 4143.227
$ws.Range("L61").Value = 7501.3335
$ws.Range("M61").Value = -3931.227
$ws.Range("N61").Value = -7925.3335

$ws.Range("H132").Value = 5093.4326
$ws.Range("I132").Value = 1558.625
$ws.Range("J132").Value = 7786.619
$ws.Range("K132").Value = 4675.875
$ws.Range("L132").Value = 23359.857
$ws.Range("M132").Value = -2145.875
$ws.Range("N132").Value = -28419.857

$ws.Range("H136").Value = 5118.161
$ws.Range("I136").Value = 4143.227
$ws.Range("J136").Value = 7501.3335
$ws.Range("K136").Value = 12429.681
$ws.Range("L136").Value = 22504.0005
$ws.Range("M136").Value = -9879.681
$ws.Range("N136").Value = -27604.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 187.86957
$ws.Range("J80").Value = 200
$ws.Range("L80").Value = 200
$ws.Range("N80").Value = -2196

$ws.Range("H83").Value = 187.86957
$ws.Range("J83").Value = 200
$ws.Range("L83").Value = 1000
$ws.Range("N83").Value = -10984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 22073.334
$ws.Range("I93").Value = 22073.334
$ws.Range("K93").Value = 22073.334
$ws.Range("M93").Value = -20201.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 900
$ws.Range("I121").Value = 166.66667
$ws.Range("J121").Value = 2000
$ws.Range("K121").Value = 500.00001
$ws.Range("L121").Value = 6000
$ws.Range("M121").Value = 809.99999
$ws.Range("N121").Value = -8620

$ws.Range("H131").Value = 1088.8684
$ws.Range("I131").Value = 1258.3334
$ws.Range("J131").Value = 978.34784
$ws.Range("K131").Value = 3775.0002
$ws.Range("L131").Value = 2935.04352
$ws.Range("M131").Value = 1264.9998
$ws.Range("N131").Value = -13015.04352

$ws.Range("H134").Value = 4003.5405
$ws.Range("I134").Value = 4167
$ws.Range("J134").Value = 3864.6
$ws.Range("K134").Value = 12501
$ws.Range("L134").Value = 11593.8
$ws.Range("M134").Value = -7431
$ws.Range("N134").Value = -21733.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 250007500
$ws.Range("I11").Value = 1000000000
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 1000000000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -999999861
$ws.Range("N11").Value = -10278

$ws.Range("H12").Value = 10003.667
$ws.Range("J12").Value = 10003.667
$ws.Range("L12").Value = 10003.667
$ws.Range("N12").Value = -10283.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1833.4445
$ws.Range("I16").Value = 2283.5
$ws.Range("J16").Value = 933.3333
$ws.Range("K16").Value = 2283.5
$ws.Range("L16").Value = 933.3333
$ws.Range("M16").Value = -2113.5
$ws.Range("N16").Value = -1273.3333

$ws.Range("H82").Value = 1330.3077
$ws.Range("I82").Value = 1207.8334
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1207.8334
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -846.8334
$ws.Range("N82").Value = -3522

$ws.Range("H85").Value = 1330.3077
$ws.Range("I85").Value = 1207.8334
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1207.8334
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = 40.16660000000002
$ws.Range("N85").Value = -5296

$ws.Range("H94").Value = 48543.332
$ws.Range("J94").Value = 48543.332
$ws.Range("L94").Value = 48543.332
$ws.Range("N94").Value = -49895.332

$ws.Range("H132").Value = 3087.743
$ws.Range("I132").Value = 2161.3914
$ws.Range("J132").Value = 4863.25
$ws.Range("K132").Value = 6484.174199999999
$ws.Range("L132").Value = 14589.75
$ws.Range("M132").Value = -3954.174199999999
$ws.Range("N132").Value = -19649.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1387
$ws.Range("N7").ClearContents()

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H123").Value = 37375
$ws.Range("J123").Value = 37375
$ws.Range("L123").Value = 37375
$ws.Range("N123").Value = -47175

$ws.Range("H125").Value = 65000
$ws.Range("J125").Value = 65000
$ws.Range("L125").Value = 65000
$ws.Range("N125").Value = -74840
